# Configure announcement via config file
#
# Rewrites the closing block of the fax template so that the free-text
# "Würzburger ... Workshop" announcement paragraph (and the {Signer}
# paragraph around it) are replaced by template placeholders driven by
# a config file: {Announcement}, {#HasCommentOrAnnouncement} /
# {#HasComment} ... {Comment} ... {/HasComment}{/HasCommentOrAnnouncement}.
# Also strips a bunch of now-unused legacy pPr direct formatting
# (widowControl / tabs / autoSpaceDE / autoSpaceDN / adjustRightInd /
# spacing) from the surrounding paragraphs.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- paragraph 137 (empty paragraph right before "Mit freundlichen Grüßen") ---
# Drop widowControl/tabs/autoSpaceDE/autoSpaceDN/adjustRightInd/spacing, keep rPr.
$p137 = $d.Paragraphs.Item(137)
$xml137 = @"
<w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr></w:p>
"@
$p137.Range.InsertXML($xml137)

# --- paragraph 138 ("Mit freundlichen Grüßen" / INCLUDETEXT Signature.docx field) ---
# Same pPr simplification; runs (the field + text) are unchanged.
$p138 = $d.Paragraphs.Item(138)
$p138InstrText = ' INCLUDETEXT  "D:' + '\\' + 'Development' + '\\' + 'NRZMHiDB' + '\\' + 'HaemophilusWeb' + '\\' + 'ReportTemplates' + '\\' + 'includes' + '\\' + 'Signature.docx" '
$xml138 = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:instrText xml:space="preserve">' + $p138InstrText + '</w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Mit freundlichen Grüßen</w:t></w:r></w:p>'
$p138.Range.InsertXML($xml138)

# --- paragraph 139 (empty paragraph) ---
# pPr simplified like the others; it now carries a single run with a line break.
$p139 = $d.Paragraphs.Item(139)
$xml139 = @"
<w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:br/></w:r></w:p>
"@
$p139.Range.InsertXML($xml139)

# --- paragraph 140 ({Signer} paragraph) ---
# Loses the bottom border / tabs / spacing / bold+sz18 default pPr rPr and the
# leading/trailing <w:br/> runs; gains the {#HasCommentOrAnnouncement} tag run.
$p140 = $d.Paragraphs.Item(140)
$xml140 = @"
<w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>{Signer}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr><w:t>{#HasCommentOrAnnouncement}</w:t></w:r></w:p>
"@
$p140.Range.InsertXML($xml140)

# --- paragraph 141 ("10. Würzburger Meningokokken-Workshop..." paragraph) ---
# Becomes the {Announcement} placeholder paragraph, with a top border and
# new spacing, plus the {#HasComment} tag run.
$p141 = $d.Paragraphs.Item(141)
$xml141 = @"
<w:p $wNs><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr><w:spacing w:before="120"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>{Announcement}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>{#HasComment}</w:t></w:r></w:p>
"@
$p141.Range.InsertXML($xml141)

# --- brand-new "Kommentar" paragraph, inserted right after paragraph 141 ---
$p141 = $d.Paragraphs.Item(141)
$p141.Range.InsertParagraphAfter()
$pKommentar = $d.Paragraphs.Item(142)
$xmlKommentar = @"
<w:p $wNs><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>Kommentar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>: {Comment}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>{/HasComment}{/HasCommentOrAnnouncement}</w:t></w:r></w:p>
"@
$pKommentar.Range.InsertXML($xmlKommentar)

# --- paragraph after Kommentar (was 142, the fldChar "end" paragraph, now 143) ---
# Same pPr simplification as the earlier paragraphs; the field-end run is unchanged.
$pEnd = $d.Paragraphs.Item(143)
$xmlEnd = @"
<w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>
"@
$pEnd.Range.InsertXML($xmlEnd)

Write-Host "Done."
